$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Edit 1: insert 16 new paragraphs right after the (empty) paragraph
# that follows "Statement: SPO Intersection (of the three sets)" and
# right before the bold "Models: Sets, Individuals, Mappings" heading.
# ------------------------------------------------------------------

$anchorSearch = $d.Content
$anchorSearch.Find.Execute("Statement: SPO Intersection (of the three sets)") | Out-Null

if ($anchorSearch.Find.Found) {
    $anchorPara = $anchorSearch.Paragraphs(1)
    # the blank paragraph immediately following the "Statement: ..." line
    $insertIndex = $anchorPara.Index + 1

    $newParaTexts = @(
        "Sets / Individuals Mappings:",
        "",
        "IDs: metaclass, class, instance, context, role, occurrence, previous, next ID roles relations for Model Set Contexts.",
        "",
        "Augmentations / Transforms: Model / Domains functional mappings. Order. Dimensions. Axes. Flows. Hierarchies. Inference / Population.",
        "",
        "Levels: Augmentations. Mappings.",
        "Levels: Resource, Kind, Statement.",
        "Levels: Reify Statement as Kind, Kind as Resource, Resource as Statement.",
        "Levels: Reify Resource as Kind, Kind as Statement, Statement as Resource.",
        "",
        "Sets / Individuals Mappings:",
        "Levels (layer statements) shifts (quads matrix). CSPO roles:",
        "",
        "(Dimension, Resource, Kind, Statement);",
        ""
    )

    foreach ($newText in $newParaTexts) {
        $d.Paragraphs($insertIndex).Range.InsertParagraphAfter()
        $insertIndex = $insertIndex + 1
        if ($newText -ne "") {
            $d.Paragraphs($insertIndex).Range.Text = $newText
        }
    }
}

# ------------------------------------------------------------------
# Edit 2: drop the trailing empty run in the "Diagrams (TO DO):"
# paragraph (it now has exactly one run with the text).
# ------------------------------------------------------------------

$diagSearch = $d.Content
$diagSearch.Find.Execute("Diagrams (TO DO):") | Out-Null

if ($diagSearch.Find.Found) {
    $diagPara = $diagSearch.Paragraphs(1)
    $diagRange = $diagPara.Range
    $cleanRange = $d.Range($diagRange.Start, $diagRange.End)
    $cleanRange.Text = "Diagrams (TO DO):"
}
